$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Swap country labels that changed rank order (sorted by Casos
#    totales, column B, descending). The underlying numeric data for
#    each row stays attached to its row index, only the country name
#    (and that row's own updated stats) move.
# -----------------------------------------------------------------

# Row 60 / 61 : Luxemburgo <-> Kazajistan
# Use a temporary placeholder on one side first so the two labels are
# never simultaneously identical to an existing shared string, which
# keeps the swap clean.
$ws.Range("A61").Value = "__TMP_SWAP_1__"
$ws.Range("A60").Value = "Kazajistan"
$ws.Range("A61").Value = "Luxemburgo"

# Row 176 / 177 : San Martin (Parte Francesa) <-> Malaui
$ws.Range("A177").Value = "__TMP_SWAP_2__"
$ws.Range("A176").Value = "Malaui"
$ws.Range("A177").Value = "San Martin (Parte Francesa)"

# -----------------------------------------------------------------
# 2) Update the numeric statistics (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# -----------------------------------------------------------------

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1142688
$ws.Range("C4").Value = 11658
$ws.Range("D4").Value = 162107
$ws.Range("E4").Value = 913961
$ws.Range("G4").Value = 867
$ws.Range("H4").Value = 66620

# Alemania (row 9)
$ws.Range("B9").Value = 164478
$ws.Range("C9").Value = 401
$ws.Range("E9").Value = 28742

# Irlanda (row 25)
$ws.Range("B25").Value = 21176
$ws.Range("C25").Value = 343
$ws.Range("E25").Value = 6504
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 1286

# Kazajistan (row 60, updated data, now ranked above Luxemburgo)
$ws.Range("B60").Value = 3857
$ws.Range("C60").Value = 260
$ws.Range("D60").Value = 985
$ws.Range("E60").Value = 2847
$ws.Range("F60").Value = 41
$ws.Range("H60").Value = 25

# Luxemburgo (row 61, data unchanged but now below Kazajistan)
$ws.Range("B61").Value = 3812
$ws.Range("C61").Value = 10
$ws.Range("D61").Value = 3318
$ws.Range("E61").Value = 402
$ws.Range("F61").Value = 23
$ws.Range("H61").Value = 92

# Sri Lanka (row 103)
$ws.Range("B103").Value = 702
$ws.Range("C103").Value = 12
$ws.Range("E103").Value = 523

# Montenegro (row 127)
$ws.Range("D127").Value = 245
$ws.Range("E127").Value = 70

# Malaui (row 176, updated data, now ranked above San Martin)
$ws.Range("C176").Value = 1
$ws.Range("D176").Value = 9
$ws.Range("E176").Value = 26
$ws.Range("F176").Value = 1

# San Martin (Parte Francesa) (row 177, data unchanged but now below Malaui)
$ws.Range("B177").Value = 38
$ws.Range("D177").Value = 27
$ws.Range("E177").Value = 8
$ws.Range("F177").Value = 3
